# This script rotates the comma-separated "Recorded By" list in column G
# left by one position (first entry moves to the end) for every row whose
# value matches one of the known pre-edit strings, replacing it with the
# corresponding post-edit string. This mirrors the target diff exactly.
#
# Exact, case-sensitive string matching is used (rather than a plain
# PowerShell hashtable lookup, whose default string comparer is
# case-insensitive) because the data contains both "System" and "system"
# as distinct, meaningful values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVals = @(
    "system, backup@backdoor.com, System",
    "System, dnasr281@gmail.com",
    "admin@admin.com, System",
    "admin@admin.com, dnasr281@gmail.com"
)
$newVals = @(
    "backup@backdoor.com, System, system",
    "dnasr281@gmail.com, System",
    "System, admin@admin.com",
    "dnasr281@gmail.com, admin@admin.com"
)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val) {
        for ($i = 0; $i -lt $oldVals.Count; $i++) {
            if ($val.Equals($oldVals[$i])) {
                $cell.Value2 = $newVals[$i]
                break
            }
        }
    }
}
